# Update the "register" sheet data: replace the sample registration rows
# with a new set of names (firstname/lastname columns), leaving telephone,
# password, and subscribe columns as-is except where noted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("register")

# firstname column updates
$ws.Range("A2").Value = "amir"
$ws.Range("A3").Value = "iravati"
$ws.Range("A4").Value = "anu"

# lastname column updates
$ws.Range("B2").Value = "tester"
$ws.Range("B3").Value = "tester"
$ws.Range("B4").Value = "tester"

# password column update for row 4
$ws.Range("D4").Value = "anu34"

# Move the active selection to A2 (was C5)
$ws.Range("A2").Select()
